$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column values stay as literal text (avoid numeric auto-conversion)
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D31", "D33", "D34", "D35", "D38", "D40", "D41", "D43", "D44", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.549.14'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '3.904.01'
$ws.Range('E3').Value = '  +2.63%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '602.52'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = '165.69'
$ws.Range('E6').Value = '  +1.25%  '
$ws.Range('D7').Value = '3.900.40'
$ws.Range('E7').Value = '  +2.57%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -1.31%  '
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('D11').Value = '6.41'
$ws.Range('E11').Value = '  +1.60%  '
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').Value = '0.0000256'
$ws.Range('E13').Value = '  +4.00%  '
$ws.Range('D14').Value = '37.32'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').Value = '4.556.45'
$ws.Range('E15').Value = '  +2.62%  '
$ws.Range('D16').Value = '3.888.59'
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('D17').Value = '68.646.17'
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '7.45'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = '17.07'
$ws.Range('E19').Value = '  -1.40%  '
$ws.Range('E20').Value = '  -2.38%  '
$ws.Range('D21').Value = '11.01'
$ws.Range('E21').Value = '  -2.40%  '
$ws.Range('D22').Value = '486.87'
$ws.Range('E22').Value = '  -0.57%  '
$ws.Range('D23').Value = '0.724'
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').Value = '0.0000168'
$ws.Range('E24').Value = '  +10.78%  '
$ws.Range('D25').Value = '84.44'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('D27').Value = '12.05'
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('D28').Value = '10.08'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  -1.08%  '
$ws.Range('D31').Value = '4.056.99'
$ws.Range('E31').Value = '  +2.69%  '
$ws.Range('E32').Value = '  -0.88%  '
$ws.Range('D33').Value = '7.74'
$ws.Range('E33').Value = '  -3.95%  '
$ws.Range('D34').Value = '31.87'
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('D35').Value = '3.857.15'
$ws.Range('E35').Value = '  +2.82%  '
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('E37').Value = '  +2.29%  '
$ws.Range('D38').Value = '5.93'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('E39').Value = '  -2.35%  '
$ws.Range('D40').Value = '3.18'
$ws.Range('E40').Value = '  +4.67%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  -2.51%  '
$ws.Range('D43').Value = '428.98'
$ws.Range('E43').Value = '  +1.47%  '
$ws.Range('D44').Value = '48.39'
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('D46').Value = '8.50'
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '142.62'
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '26.28'
$ws.Range('E49').Value = '  +7.65%  '
$ws.Range('D50').Value = '2.807.28'
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('D51').Value = '0.0352'
$ws.Range('E51').Value = '  +0.38%  '
